# Applies the cell-level edits described by the target diff:
# small "want to go" (F) count increments across all four sheets, plus a
# like-for-like row-content shift in the 全部类型 (all-types) rollup sheet
# so it catches up with rows already present in 展览/演出.
$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 455

$wsExhibit.Range("F8").Value = 1197

$wsExhibit.Range("F10").Value = 1305

$wsExhibit.Range("F11").Value = 879

$wsExhibit.Range("F12").Value = 686

$wsExhibit.Range("F13").Value = 185

$wsExhibit.Range("F14").Value = 508

$wsExhibit.Range("F18").Value = 2924

$wsExhibit.Range("F19").Value = 2621

$wsExhibit.Range("F24").Value = 229

$wsExhibit.Range("F31").Value = 310

# ---- Sheet: 演出 ----
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F26").Value = 3933

$wsShow.Range("F33").Value = 164

# ---- Sheet: 本地生活 ----
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F5").Value = 2456

$wsLocal.Range("F6").Value = 1043

$wsLocal.Range("F9").Value = 1322

# ---- Sheet: 全部类型 ----
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 2456

$wsAll.Range("F6").Value = 1043

$wsAll.Range("F7").Value = 1322

$wsAll.Range("F11").Value = 455

$wsAll.Range("F14").Value = 1197

$wsAll.Range("F16").Value = 879

$wsAll.Range("F17").Value = 686

$wsAll.Range("C20").Value = "上海·  第五十三届妖漫动漫游戏展"
$wsAll.Range("D20").Value = "漕溪北路339号百脑汇4楼 百脑汇"
$wsAll.Range("E20").Value = "2024.02.24 10:00-02.25 17:00"
$wsAll.Range("F20").Value = 185
$wsAll.Range("G20").Value = 80
$wsAll.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=78657"
$wsAll.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202402/MVRgtEd91707208495928.jpeg"

$wsAll.Range("C21").Value = "上海·SISPmini动漫游戏嘉年华"
$wsAll.Range("D21").Value = "剑川路1000号 龙湖上海闵行天街"
$wsAll.Range("E21").Value = "2024.02.24 13:00-02.25 19:00"
$wsAll.Range("F21").Value = 508
$wsAll.Range("G21").Value = 48
$wsAll.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=79046"
$wsAll.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202401/jzovdppq1706166165502.jpeg"

$wsAll.Range("C22").Value = "上海·《哈利的魔法世界》动漫视听音乐会"
$wsAll.Range("D22").Value = "都市路4889号（莘庄地铁站南广场） 上海保利城市剧院"
$wsAll.Range("E22").Value = "2024.02.24 14:30-02.24 16:00"
$wsAll.Range("F22").Value = 18
$wsAll.Range("G22").Value = 158
$wsAll.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=80639"
$wsAll.Range("I22").Value = "//i2.hdslb.com/bfs/openplatform/202401/4PieCC9N1706261750579.jpeg"

$wsAll.Range("C23").Value = "上海·原X铁X崩only"
$wsAll.Range("D23").Value = "澳门路168号 月星国际家居"
$wsAll.Range("E23").Value = "2024.02.24 10:30-02.24 16:30"
$wsAll.Range("F23").Value = 173
$wsAll.Range("G23").Value = 60
$wsAll.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=81446"
$wsAll.Range("I23").Value = "//i2.hdslb.com/bfs/openplatform/202401/IIePRulM1706248855263.jpeg"

$wsAll.Range("C24").Value = "上海·原神×崩坏×星铁only旅行盛宴2.0"
$wsAll.Range("D24").Value = "西藏南路1号 上海大世界"
$wsAll.Range("E24").Value = "2024.02.24 10:00-02.25 17:00"
$wsAll.Range("F24").Value = 2924
$wsAll.Range("G24").Value = 65
$wsAll.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=81276"
$wsAll.Range("I24").Value = "//i2.hdslb.com/bfs/openplatform/202401/82hU3z8m1706155835021.png"

$wsAll.Range("C25").Value = "上海·第三届燃梦BACG PRO游戏动漫展-原X铁X崩同好交流"
$wsAll.Range("D25").Value = "盈浦街道淀山浦社区淀山湖大道851号青浦万达茂F3 万达汽车乐园(青浦万达茂店)"
$wsAll.Range("E25").Value = "2024.02.24 11:00-02.25 16:30"
$wsAll.Range("F25").Value = 2621
$wsAll.Range("G25").Value = 65.8
$wsAll.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=77754"
$wsAll.Range("I25").Value = "//i1.hdslb.com/bfs/openplatform/202312/7eGZETK91701943985671.jpeg"

$wsAll.Range("F29").Value = 229

$wsAll.Range("F38").Value = 310

$wsAll.Range("C46").Value = "上海·首届Redamancy动漫游戏嘉年华"
$wsAll.Range("D46").Value = "中山北路3300号4楼L4001号 环球港上海世嘉都市乐园"
$wsAll.Range("E46").Value = "2024.03.30 10:00-03.31 17:00"
$wsAll.Range("F46").Value = 1096
$wsAll.Range("G46").Value = 60
$wsAll.Range("H46").Value = "https://show.bilibili.com/platform/detail.html?id=81772"
$wsAll.Range("I46").Value = "//i2.hdslb.com/bfs/openplatform/202402/XKf9RSFB1707127784856.jpeg"

$wsAll.Range("B47").NumberFormat = "@"
$wsAll.Range("B47").Value = "2024-04-13"
$wsAll.Range("B47").ClearFormats()
$wsAll.Range("C47").Value = "上海·《四月是你的谎言》——“公生”与“薰”的钢琴小提琴唯美经典音乐集"
$wsAll.Range("D47").Value = "丁香路425号 上海东方艺术中心"
$wsAll.Range("E47").Value = "2024.04.13 19:30-04.13 21:30"
$wsAll.Range("F47").Value = 198
$wsAll.Range("G47").Value = 80
$wsAll.Range("H47").Value = "https://show.bilibili.com/platform/detail.html?id=78667"
$wsAll.Range("I47").Value = "//i1.hdslb.com/bfs/openplatform/202311/bTP7w6GD1700130122940.jpeg"

$wsAll.Range("B48").NumberFormat = "@"
$wsAll.Range("B48").Value = "2024-04-20"
$wsAll.Range("B48").ClearFormats()
$wsAll.Range("C48").Value = "上海· 茅原实里动漫交响音乐会"
$wsAll.Range("D48").Value = "东大名路889号 友邦大剧院"
$wsAll.Range("E48").Value = "2024.04.20 19:30-04.20 21:00"
$wsAll.Range("F48").Value = 51
$wsAll.Range("G48").Value = 280
$wsAll.Range("H48").Value = "https://show.bilibili.com/platform/detail.html?id=81703"
$wsAll.Range("I48").Value = "//i2.hdslb.com/bfs/openplatform/202402/yiVaqJVK1707016321221.jpeg"

$wsAll.Range("F49").Value = 164
